$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the old last data row (row 14) -- table shrinks from 9 to 8 data rows
$ws.Rows.Item(14).Delete()

# Column C needs to widen (class name is now longer: "Testing class")
$ws.Columns.Item(3).ColumnWidth = 15.6

# Row 6: was Michael Johnson / Laskar Pelangi entry -> becomes "siswa satus" / "Testing class" entry
$ws.Range("B6").Value = "siswa satus"
$ws.Range("C6").Value = "Testing class"
$ws.Range("D6").Value = "-"
$ws.Range("E6").Value = "Buku Paket Bahasa Ingris"
$ws.Range("F6").Value = "978-0-393-04002-9"
$ws.Range("G6").Value = 15
$ws.Range("H6").Value = "26-02-2025"
$ws.Range("I6").Value = "25-02-2025"
$ws.Range("J6").Value = "Terlambat"

# Row 7: was "siswa satu" entry -> becomes Michael Johnson / Laskar Pelangi entry
$ws.Range("B7").Value = "Michael Johnson"
$ws.Range("C7").Value = "XI-A"
$ws.Range("D7").Value = 2147483647
$ws.Range("E7").Value = "Laskar Pelangi"
$ws.Range("F7").Value = "978-3-16-148410-0"
$ws.Range("G7").Value = 43
$ws.Range("H7").Value = "26-02-2025"
# "06-03-2025" looks like a valid M-D-Y date, so force text formatting first
# to avoid Excel silently converting it to a date serial number.
$ws.Range("I7").NumberFormat = "@"
$ws.Range("I7").Value = "06-03-2025"
$ws.Range("H7").Copy()
$ws.Range("I7").PasteSpecial(-4122)
$ws.Range("J7").Value = "Terlambat"

# Row 8: unchanged

# Row 9: Jumlah Pinjaman changes
$ws.Range("G9").Value = 32

# Row 10: Jumlah Pinjaman + Status change
$ws.Range("G10").Value = 2
$ws.Range("J10").Value = "Terlambat"

# Row 11: Status changes
$ws.Range("J11").Value = "Dikembalikan"

# Row 12: Status changes
$ws.Range("J12").Value = "Diperpanjang"

# Row 13: was Michael Johnson entry -> becomes "siswa satus" / "Testing class" entry (merged with old row 14's note)
$ws.Range("B13").Value = "siswa satus"
$ws.Range("C13").Value = "Testing class"
$ws.Range("D13").Value = "-"
$ws.Range("G13").Value = 21
$ws.Range("H13").Value = "24-02-2025"
# "08-03-2025" looks like a valid M-D-Y date too, same fix as I7 above.
$ws.Range("I13").NumberFormat = "@"
$ws.Range("I13").Value = "08-03-2025"
$ws.Range("H13").Copy()
$ws.Range("I13").PasteSpecial(-4122)
$ws.Range("J13").Value = "Diperpanjang"
$ws.Range("K13").Value = "jangan sampai telat"

# Keep the header/data selection in sync with the now-shorter table
$ws.Range("A5:M13").Select() | Out-Null
